$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Edit 1: "What is indirection operator?" -- split the " operator" run into
# two runs: " " and "operator".
# ---------------------------------------------------------------------------
$r = $d.Content.Duplicate
$r.Find.Execute("What is indirection operator?", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Text = ""
$ins = $d.Range($r.Start, $r.Start)
$ins.InsertXML("<w:p $wNs><w:r><w:t>What is indirection</w:t></w:r><w:r><w:t xml:space=`"preserve`"> </w:t></w:r><w:r><w:t>operator</w:t></w:r><w:r><w:t>?</w:t></w:r></w:p>")

Write-Output "edit1 done"

# ---------------------------------------------------------------------------
# Edit 2: delete the four questions about pointer levels / sentinel value /
# void pointer / far & near pointer (the list item right after them, about
# structure vs union, keeps its existing formatting).
# ---------------------------------------------------------------------------
$s1 = $d.Content.Duplicate
$s1.Find.Execute("How many levels of pointer can you have?", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s2 = $d.Content.Duplicate
$s2.Find.Execute("What is the difference between structure and union?", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$toDelete = $d.Range($s1.Start, $s2.Start)
$toDelete.Delete()

Write-Output "edit2 done"
